# FIX: Se corrige puntaje y se añaden resultados
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    1  = 8
    2  = 6
    3  = 7
    4  = 4
    5  = 4
    7  = 2
    8  = 2
    9  = 1
    14 = 0
    16 = 17
    17 = 5
    18 = 3
    19 = 2
    20 = 4
    21 = 5
    23 = 1
    35 = 1
    38 = 1
    41 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
